# Rebuild the "保險" (insurance, sheet index 5) and "債務" (debt, sheet index 6)
# sheets with the corrected/normalized schema (company/species/insurance columns
# for sheet5; species/debtor/debt columns for sheet6), matching the data
# pipeline's "insurance, claim, debt, investment done" pass.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet5 : 保險 (insurance)
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Cells.Clear()

$idx5 = @(96, 97, 98, 99, 100, 101, 102)
$company5 = @("南山人壽", "南山人壽", "台灣人壽", "大都會國際人壽", "國泰人壽", "國泰人壽", "幸福人壽")
$species5 = @("新康祥终身壽險", "財星高照變額萬能壽險", "富利人生終身壽險B型", "金多多保險", "萬代福211", "雙星還本", "大吉大利终身壽險")
$owner5 = @("簡東明", "簡東明", "戴錦花", "戴錦花", "戴錦花", "戴錦花", "戴錦花")

# Header row
$ws5.Cells.Item(1, 2).Value = "company"
$ws5.Cells.Item(1, 3).Value = "name"
$ws5.Cells.Item(1, 4).Value = "owner"
$ws5.Cells.Item(1, 5).Value = "property_category"
$ws5.Cells.Item(1, 6).Value = "category"
$ws5.Cells.Item(1, 7).Value = "date"
$ws5.Cells.Item(1, 8).Value = "legislator_name"
$ws5.Cells.Item(1, 9).Value = "legislator_id"
$ws5.Cells.Item(1, 10).Value = "source_file"
$ws5.Cells.Item(1, 11).Value = "index"

# Data, written column by column (A, then B, then C, ...)
for ($r = 0; $r -lt $idx5.Length; $r++) {
    $ws5.Cells.Item($r + 2, 1).Value = $idx5[$r]
}
for ($r = 0; $r -lt $company5.Length; $r++) {
    $ws5.Cells.Item($r + 2, 2).Value = $company5[$r]
}
for ($r = 0; $r -lt $species5.Length; $r++) {
    $ws5.Cells.Item($r + 2, 3).Value = $species5[$r]
}
for ($r = 0; $r -lt $owner5.Length; $r++) {
    $ws5.Cells.Item($r + 2, 4).Value = $owner5[$r]
}
for ($r = 0; $r -lt $idx5.Length; $r++) {
    $ws5.Cells.Item($r + 2, 5).Value = "insurance"
}
for ($r = 0; $r -lt $idx5.Length; $r++) {
    $ws5.Cells.Item($r + 2, 6).Value = "normal"
}
for ($r = 0; $r -lt $idx5.Length; $r++) {
    $ws5.Cells.Item($r + 2, 7).Value = "2011-12-30"
}
for ($r = 0; $r -lt $idx5.Length; $r++) {
    $ws5.Cells.Item($r + 2, 8).Value = "簡東明"
}
for ($r = 0; $r -lt $idx5.Length; $r++) {
    $ws5.Cells.Item($r + 2, 9).Value = 1717
}
for ($r = 0; $r -lt $idx5.Length; $r++) {
    $ws5.Cells.Item($r + 2, 10).Value = "tmp3d8a1"
}
for ($r = 0; $r -lt $idx5.Length; $r++) {
    $ws5.Cells.Item($r + 2, 11).Value = $idx5[$r]
}

# ---------------------------------------------------------------------
# Sheet6 : 債務 (debt)
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)
$ws6.Cells.Clear()

# Header row
$ws6.Cells.Item(1, 2).Value = "species"
$ws6.Cells.Item(1, 3).Value = "debtor"
$ws6.Cells.Item(1, 4).Value = "owner"
$ws6.Cells.Item(1, 5).Value = "total"
$ws6.Cells.Item(1, 6).Value = "register_date"
$ws6.Cells.Item(1, 7).Value = "register_reason"
$ws6.Cells.Item(1, 8).Value = "property_category"
$ws6.Cells.Item(1, 9).Value = "category"
$ws6.Cells.Item(1, 10).Value = "date"
$ws6.Cells.Item(1, 11).Value = "legislator_name"
$ws6.Cells.Item(1, 12).Value = "legislator_id"
$ws6.Cells.Item(1, 13).Value = "source_file"
$ws6.Cells.Item(1, 14).Value = "index"

# Single data row
$ws6.Cells.Item(2, 1).Value = 112
$ws6.Cells.Item(2, 2).Value = "房屋貸款"
$ws6.Cells.Item(2, 3).Value = "戴錦花"
$ws6.Cells.Item(2, 4).Value = "永豐銀行屏東縣屏東市復興北路"
$ws6.Cells.Item(2, 5).Value = 1120295
$ws6.Cells.Item(2, 6).Value = "92年07月01日"
$ws6.Cells.Item(2, 7).Value = "購屋"
$ws6.Cells.Item(2, 8).Value = "debt"
$ws6.Cells.Item(2, 9).Value = "normal"
$ws6.Cells.Item(2, 10).Value = "2011-12-30"
$ws6.Cells.Item(2, 11).Value = "簡東明"
$ws6.Cells.Item(2, 12).Value = 1717
$ws6.Cells.Item(2, 13).Value = "tmp3d8a1"
$ws6.Cells.Item(2, 14).Value = 112
